$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "No Bait Consumption" translation block (rows 111-113) ---
# Mirrors the existing title/toggle/tooltip row pattern used for every other
# feature in this sheet (e.g. rows 108-110 for "No Stamina Cost").

$ws.Range("A111").Value = "title22"
$ws.Range("B111").Value = "No Bait Consumption"
$ws.Range("C111").Value = "餌消費なし"
$ws.Range("D111").Value = "不消耗鱼饵"

$ws.Range("A112").Value = "toggle56"
$ws.Range("B112").Value = "Enable No Bait Consumption"
$ws.Range("C112").Value = "餌消費なしを有効化"
$ws.Range("D112").Value = "启用不消耗鱼饵"

$ws.Range("A113").Value = "tooltip22"
$ws.Range("B113").Value = "Enable or disable no bait consumption while fishing."
$ws.Range("C113").Value = "釣り中の餌消費なしを有効または無効にします。"
$ws.Range("D113").Value = "启用或禁用钓鱼时不消耗鱼饵。"

# The Japanese/Chinese columns need the CJK-capable font used throughout the
# rest of the sheet (columns A/B already pick up the correct default font).
foreach ($addr in @("C111", "D111", "C112", "D112", "C113", "D113")) {
    $ws.Range($addr).Font.Name = "Noto Sans SC"
}

# --- Normalise the font variant on the existing rows touched by this edit ---
# Rows 105 and 108-110 had their C/D cells re-saved with the same CJK font.
foreach ($addr in @("C105", "C108", "D108", "C109", "D109", "C110", "D110")) {
    $ws.Range($addr).Font.Name = "Noto Sans SC"
}

# --- Keep the view/selection state consistent with the now-larger sheet ---
$ws.Range("C102").Select()
